$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 data (CW3M Baseline 2010-18 C371+ run)
$ws.Range("A6").Value = "CW3M"
$ws.Range("B6").Value = "Baseline 2010-18 C371+"
$ws.Range("C6").Value = "2010-18"

$ws.Range("D6").Value = 675.83090190000007
$ws.Range("D6").NumberFormat = "0.00"

$ws.Range("E6").Value = 2040.5741823000001
$ws.Range("E6").NumberFormat = "0.00"
$ws.Range("E6").Interior.Color = 65535

$ws.Range("F6").Value = 5.7945578999999992
$ws.Range("F6").NumberFormat = "0.00"
$ws.Range("F6").Interior.Color = 65535

$ws.Range("G6").Value = 232.20442180000001
$ws.Range("G6").NumberFormat = "0.00"

$ws.Range("H6").Value = 0
$ws.Range("H6").NumberFormat = "0.00"

$ws.Range("I6").Value = 6.2224744999999997
$ws.Range("I6").NumberFormat = "0.00"

$ws.Range("J6").Value = 0
$ws.Range("J6").NumberFormat = "0.00"

$ws.Range("K6").Value = 549.56830450000007
$ws.Range("K6").NumberFormat = "0.00"
$ws.Range("K6").Interior.Color = 65535

$ws.Range("L6").Value = 86.997628399999996
$ws.Range("L6").NumberFormat = "0.00"
$ws.Range("L6").Interior.Color = 65535

$ws.Range("M6").Value = 1652.2535766000001
$ws.Range("M6").NumberFormat = "0.00"

$ws.Range("N6").Value = 668.55730879999999
$ws.Range("N6").NumberFormat = "0.00"
$ws.Range("N6").Interior.Color = 65535

$ws.Range("O6").Value = 15727.597461100002
$ws.Range("O6").NumberFormat = "0"
$ws.Range("O6").Interior.Color = 65535

$ws.Range("P6").Value = 2215.5502928999999
$ws.Range("P6").NumberFormat = "0"

$ws.Range("Q6").Value = -3.2497201000000002
$ws.Range("Q6").NumberFormat = "0.00"

$ws.Range("R6").Value = [double]"-1.3625999999999998E-3"
$ws.Range("R6").NumberFormat = "0.000000"

$ws.Range("B6").Select()
